$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the "estado de cuenta" detail rows (16-21) with the new data set.
# Row 16: KARINA DEL CARMEN VILLALBA BORJA - periodo 1803 - valor mora 36000
$ws.Range("C16").Value = "30873862"
$ws.Range("D16").Value = "KARINA DEL CARMEN VILLALBA BORJA"
$ws.Range("E16").Value = "1803"
$ws.Range("F16").Value = 36000
$ws.Range("G16").Value = 900000

# Row 17: KARINA DEL CARMEN VILLALBA BORJA - periodo 1803 - valor mora 19200
$ws.Range("C17").Value = "30873862"
$ws.Range("D17").Value = "KARINA DEL CARMEN VILLALBA BORJA"
$ws.Range("E17").Value = "1803"
$ws.Range("F17").Value = 19200
$ws.Range("G17").Value = 900000

# Row 18: PAOLA PATRICIA AGUILAR VELASCO - periodo 1804 - valor mora 31249
$ws.Range("C18").Value = "1143379280"
$ws.Range("D18").Value = "PAOLA PATRICIA AGUILAR VELASCO"
$ws.Range("E18").Value = "1804"
$ws.Range("F18").Value = 31249
$ws.Range("G18").Value = 781242

# Row 19: PAOLA PATRICIA AGUILAR VELASCO - periodo 1803 - valor mora 31249
$ws.Range("C19").Value = "1143379280"
$ws.Range("D19").Value = "PAOLA PATRICIA AGUILAR VELASCO"
$ws.Range("E19").Value = "1803"
$ws.Range("F19").Value = 31249
$ws.Range("G19").Value = 781242

# Row 20: DARIANY CANO DIAZ - periodo 1804 - valor mora 30208
$ws.Range("C20").Value = "1143404861"
$ws.Range("D20").Value = "DARIANY CANO DIAZ"
$ws.Range("E20").Value = "1804"
$ws.Range("F20").Value = 30208
$ws.Range("G20").Value = 781242

# Row 21: BIBIANA LOPEZ DIAZ - periodo 1804 - valor mora 31249
$ws.Range("C21").Value = "30656412"
$ws.Range("D21").Value = "BIBIANA LOPEZ DIAZ"
$ws.Range("E21").Value = "1804"
$ws.Range("F21").Value = 31249
$ws.Range("G21").Value = 781242
